$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 32 (item id 5484) on ALC
$ws_ALC.Range("H32").Value = 2999.5
$ws_ALC.Range("J32").Value = 2999.5
$ws_ALC.Range("L32").Value = 2999.5
$ws_ALC.Range("N32").Value = -3651.5

# Row 87 (item id 10651) on ALC
$ws_ALC.Range("H87").Value = 0
$ws_ALC.Range("J87").Value = 0
$ws_ALC.Range("L87").Value = 0
$ws_ALC.Range("N87").Value = $null

# Row 90 (item id 10651) on ALC
$ws_ALC.Range("H90").Value = 0
$ws_ALC.Range("J90").Value = 0
$ws_ALC.Range("L90").Value = 0
$ws_ALC.Range("N90").Value = $null

# Row 118 (item id 27958) on ALC
$ws_ALC.Range("H118").Value = 169.5
$ws_ALC.Range("I118").Value = 169.5
$ws_ALC.Range("K118").Value = 508.5
$ws_ALC.Range("M118").Value = 1148.5

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 13 (item id 2656) on ARM
$ws_ARM.Range("H13").Value = 5001000
$ws_ARM.Range("I13").Value = 5001000
$ws_ARM.Range("K13").Value = 5001000
$ws_ARM.Range("M13").Value = -5000856

# Row 37 (item id 3096) on ARM
$ws_ARM.Range("H37").Value = 23017
$ws_ARM.Range("I37").Value = 1034
$ws_ARM.Range("K37").Value = 1034
$ws_ARM.Range("M37").Value = -761

# Row 132 (item id 43997) on ARM
$ws_ARM.Range("H132").Value = 10833.167
$ws_ARM.Range("I132").Value = 3999.6667
$ws_ARM.Range("K132").Value = 11999.0001
$ws_ARM.Range("M132").Value = -9469.000100000001

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 11 (item id 2481) on BSM
$ws_BSM.Range("H11").Value = 996
$ws_BSM.Range("J11").Value = 996
$ws_BSM.Range("L11").Value = 996
$ws_BSM.Range("N11").Value = -1276

# Row 20 (item id 14149) on BSM
$ws_BSM.Range("H20").Value = 3325.7144
$ws_BSM.Range("I20").Value = 3325.7144
$ws_BSM.Range("K20").Value = 3325.7144
$ws_BSM.Range("M20").Value = -3078.7144

# Row 94 (item id 19939) on BSM
$ws_BSM.Range("H94").Value = 1999.8
$ws_BSM.Range("I94").Value = 1999.75
$ws_BSM.Range("K94").Value = 1999.75
$ws_BSM.Range("M94").Value = -1548.75

# Row 99 (item id 19943) on BSM
$ws_BSM.Range("H99").Value = 2737.5
$ws_BSM.Range("I99").Value = 2737.5
$ws_BSM.Range("K99").Value = 2737.5
$ws_BSM.Range("M99").Value = -1239.5

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 10 (item id 1997) on CRP
$ws_CRP.Range("H10").Value = 431.2
$ws_CRP.Range("I10").Value = 264
$ws_CRP.Range("K10").Value = 264
$ws_CRP.Range("M10").Value = -125

# Row 88 (item id 10608) on CRP
$ws_CRP.Range("H88").Value = 36150
$ws_CRP.Range("J88").Value = 36150
$ws_CRP.Range("L88").Value = 36150
$ws_CRP.Range("N88").Value = -36962

# Row 91 (item id 10608) on CRP
$ws_CRP.Range("H91").Value = 36150
$ws_CRP.Range("J91").Value = 36150
$ws_CRP.Range("L91").Value = 36150
$ws_CRP.Range("N91").Value = -38958

# Row 92 (item id 18041) on CRP
$ws_CRP.Range("H92").Value = 9189
$ws_CRP.Range("J92").Value = 9189
$ws_CRP.Range("L92").Value = 9189
$ws_CRP.Range("N92").Value = -14181

# Row 95 (item id 18192) on CRP
$ws_CRP.Range("H95").Value = 13450
$ws_CRP.Range("J95").Value = 13450
$ws_CRP.Range("L95").Value = 13450
$ws_CRP.Range("N95").Value = -18942

# Row 96 (item id 18193) on CRP
$ws_CRP.Range("H96").Value = 9972.6
$ws_CRP.Range("J96").Value = 9972.6
$ws_CRP.Range("L96").Value = 9972.6
$ws_CRP.Range("N96").Value = -15464.6

# Row 107 (item id 27689) on CRP
$ws_CRP.Range("H107").Value = 4526.2104
$ws_CRP.Range("I107").Value = 4617.5884
$ws_CRP.Range("K107").Value = 4617.5884
$ws_CRP.Range("M107").Value = -2697.5884

# Row 134 (item id 44020) on CRP
$ws_CRP.Range("H134").Value = 6577.6
$ws_CRP.Range("I134").Value = 3793.8333
$ws_CRP.Range("J134").Value = 10753.25
$ws_CRP.Range("K134").Value = 11381.4999
$ws_CRP.Range("L134").Value = 32259.75
$ws_CRP.Range("M134").Value = -8846.499899999999
$ws_CRP.Range("N134").Value = -37329.75

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 4 (item id 4650) on CUL
$ws_CUL.Range("H4").Value = 142859100
$ws_CUL.Range("I4").Value = 2287.8333
$ws_CUL.Range("K4").Value = 6863.499899999999
$ws_CUL.Range("M4").Value = -6751.499899999999

# Row 23 (item id 4858) on CUL
$ws_CUL.Range("H23").Value = 209.5
$ws_CUL.Range("I23").Value = 215.5
$ws_CUL.Range("K23").Value = 646.5
$ws_CUL.Range("M23").Value = -411.5

# Row 80 (item id 12890) on CUL
$ws_CUL.Range("H80").Value = 5034
$ws_CUL.Range("I80").Value = 2549.5
$ws_CUL.Range("K80").Value = 7648.5
$ws_CUL.Range("M80").Value = -6712.5

# Row 83 (item id 12890) on CUL
$ws_CUL.Range("H83").Value = 5034
$ws_CUL.Range("I83").Value = 2549.5
$ws_CUL.Range("K83").Value = 22945.5
$ws_CUL.Range("M83").Value = -18265.5

# Row 122 (item id 36078) on CUL
$ws_CUL.Range("H122").Value = 963.8182
$ws_CUL.Range("I122").Value = 900.25
$ws_CUL.Range("J122").Value = 1133.3334
$ws_CUL.Range("K122").Value = 8102.25
$ws_CUL.Range("L122").Value = 10200.0006
$ws_CUL.Range("M122").Value = -5652.25
$ws_CUL.Range("N122").Value = -15100.0006

# Row 131 (item id 36060) on CUL
$ws_CUL.Range("H131").Value = 14407
$ws_CUL.Range("I131").Value = 100000
$ws_CUL.Range("J131").Value = 2179.4285
$ws_CUL.Range("K131").Value = 300000
$ws_CUL.Range("L131").Value = 6538.2855
$ws_CUL.Range("M131").Value = -294960
$ws_CUL.Range("N131").Value = -16618.2855

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 92 (item id 18094) on GSM
$ws_GSM.Range("H92").Value = 9099.75
$ws_GSM.Range("J92").Value = 9099.75
$ws_GSM.Range("L92").Value = 9099.75
$ws_GSM.Range("N92").Value = -12843.75

# Row 99 (item id 19532) on GSM
$ws_GSM.Range("H99").Value = 7879.4
$ws_GSM.Range("I99").Value = 7879.4
$ws_GSM.Range("J99").Value = 0
$ws_GSM.Range("K99").Value = 7879.4
$ws_GSM.Range("L99").Value = 0
$ws_GSM.Range("M99").Value = -5633.4
$ws_GSM.Range("N99").Value = $null

# Row 126 (item id 36184) on GSM
$ws_GSM.Range("H126").Value = 8250.5
$ws_GSM.Range("I126").Value = 7667.3335
$ws_GSM.Range("J126").Value = 10000
$ws_GSM.Range("K126").Value = 23002.0005
$ws_GSM.Range("L126").Value = 30000
$ws_GSM.Range("M126").Value = -20532.0005
$ws_GSM.Range("N126").Value = -34940

# Row 134 (item id 42064) on GSM
$ws_GSM.Range("H134").Value = 0
$ws_GSM.Range("J134").Value = 0
$ws_GSM.Range("L134").Value = 0
$ws_GSM.Range("N134").Value = $null

# Row 135 (item id 42006) on GSM
$ws_GSM.Range("H135").Value = 0
$ws_GSM.Range("J135").Value = 0
$ws_GSM.Range("L135").Value = 0
$ws_GSM.Range("N135").Value = $null

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 2 (item id 2631) on LTW
$ws_LTW.Range("H2").Value = 308.81818
$ws_LTW.Range("I2").Value = 308.81818
$ws_LTW.Range("J2").Value = 0
$ws_LTW.Range("K2").Value = 308.81818
$ws_LTW.Range("L2").Value = 0
$ws_LTW.Range("M2").Value = -196.81818
$ws_LTW.Range("N2").Value = $null

# Row 22 (item id 5277) on LTW
$ws_LTW.Range("H22").Value = 2218.3635
$ws_LTW.Range("I22").Value = 925
$ws_LTW.Range("K22").Value = 925
$ws_LTW.Range("M22").Value = -630

# Row 27 (item id 5277) on LTW
$ws_LTW.Range("H27").Value = 2218.3635
$ws_LTW.Range("I27").Value = 925
$ws_LTW.Range("K27").Value = 925
$ws_LTW.Range("M27").Value = -818

# Row 93 (item id 19993) on LTW
$ws_LTW.Range("H93").Value = 2979.9
$ws_LTW.Range("I93").Value = 2949.5
$ws_LTW.Range("K93").Value = 2949.5
$ws_LTW.Range("M93").Value = -1701.5

# Row 122 (item id 36247) on LTW
$ws_LTW.Range("H122").Value = 3800
$ws_LTW.Range("J122").Value = 3800
$ws_LTW.Range("L122").Value = 11400
$ws_LTW.Range("N122").Value = -16300

# Row 135 (item id 42036) on LTW
$ws_LTW.Range("H135").Value = 32499.5
$ws_LTW.Range("J135").Value = 32499.5
$ws_LTW.Range("L135").Value = 32499.5
$ws_LTW.Range("N135").Value = -42639.5

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 2 (item id 3307) on WVR
$ws_WVR.Range("H2").Value = 4597.2
$ws_WVR.Range("I2").Value = 4597.2
$ws_WVR.Range("K2").Value = 4597.2
$ws_WVR.Range("M2").Value = -4485.2

# Row 94 (item id 18075) on WVR
$ws_WVR.Range("H94").Value = 27664.334
$ws_WVR.Range("I94").Value = 22993
$ws_WVR.Range("J94").Value = 30000
$ws_WVR.Range("K94").Value = 22993
$ws_WVR.Range("L94").Value = 30000
$ws_WVR.Range("M94").Value = -22092
$ws_WVR.Range("N94").Value = -31802
